{"js": "// The author fixed a typo in the EPICS paragraph: \"\u0438\u0441\u043f\u043e\u043b\u044c\u0437\u0443\u044e\u0438\u0441\u044f\" -> \"\u0438\u0441\u043f\u043e\u043b\u044c\u0437\u0443\u044e\u0442\u0441\u044f\"\n// (\"\u041e\u043d\u0438 \u0438\u0441\u043f\u043e\u043b\u044c\u0437\u0443\u044e\u0438\u0441\u044f \u0432 \u043d\u0430\u043f\u0438\u0441\u0430\u043d\u0438\u0438 \u0441\u0446\u0435\u043d\u0430\u0440\u0438\u0435\u0432...\" -> \"\u041e\u043d\u0438 \u0438\u0441\u043f\u043e\u043b\u044c\u0437\u0443\u044e\u0442\u0441\u044f \u0432 \u043d\u0430\u043f\u0438\u0441\u0430\u043d\u0438\u0438...\").\n// Word re-seated its hidden \"_GoBack\" (last-edit-position) bookmark from its old spot\n// (an otherwise-empty paragraph further down) to the exact point of this edit, right\n// after the inserted letter. We reproduce both effects below.\n\nconst doc = context.document;\n\n// 1) Drop the old \"_GoBack\" bookmark wherever it currently lives (Word only ever\n//    keeps one), so it doesn't linger at its old location once we re-add it below.\nconst existingGoBack = doc.getBookmarkRangeOrNullObject(\"_GoBack\");\nexistingGoBack.load(\"isNullObject\");\nawait context.sync();\nif (!existingGoBack.isNullObject) {\n  doc.deleteBookmark(\"_GoBack\");\n  await context.sync();\n}\n\n// 2) Locate the misspelled word.\nconst typoResults = doc.body.search(\"\u0438\u0441\u043f\u043e\u043b\u044c\u0437\u0443\u044e\u0438\u0441\u044f\", { matchCase: false });\ntypoResults.load(\"text\");\nawait context.sync();\nconst typo = typoResults.items[0];\n\n// 3) Within it, isolate the trailing \"\u0438\u0441\u044f\" and then just its leading \"\u0438\" - that is\n//    the single character that needs to become \"\u0442\" (\"\u0438\u0441\u043f\u043e\u043b\u044c\u0437\u0443\u044e\" + \"\u0438\u0441\u044f\" -> \"\u0438\u0441\u043f\u043e\u043b\u044c\u0437\u0443\u044e\" + \"\u0442\u0441\u044f\").\nconst tailResults = typo.search(\"\u0438\u0441\u044f\", { matchCase: false });\ntailResults.load(\"text\");\nawait context.sync();\nconst tail = tailResults.items[0];\n\nconst badCharResults = tail.search(\"\u0438\", { matchCase: false });\nbadCharResults.load(\"text\");\nawait context.sync();\nconst badChar = badCharResults.items[0];\n\nbadChar.insertText(\"\u0442\", Word.InsertLocation.replace);\nawait context.sync();\n\n// Re-load the (now corrected) range and nudge its formatting so the host keeps this\n// \"\u0442\" as its own run instead of silently re-merging it with its neighbor - matching\n// the author's document, where it stayed a distinct run next to the new bookmark.\nbadChar.load(\"font/bold\");\nawait context.sync();\nconst originalBold = badChar.font.bold;\nbadChar.font.bold = !originalBold;\nawait context.sync();\nbadChar.font.bold = originalBold;\nawait context.sync();\n\n// 4) Re-insert \"_GoBack\" as a zero-length bookmark immediately after the new \"\u0442\",\n//    i.e. right before \"\u0441\u044f ...\" - exactly where Word left the cursor after this edit.\nconst afterBadChar = badChar.getRange(\"After\");\nafterBadChar.insertBookmark(\"_GoBack\");\nawait context.sync();\n", "ps1": "# The author fixed a typo in the EPICS paragraph: \"\u0438\u0441\u043f\u043e\u043b\u044c\u0437\u0443\u044e\u0438\u0441\u044f\" -> \"\u0438\u0441\u043f\u043e\u043b\u044c\u0437\u0443\u044e\u0442\u0441\u044f\"\n# (\"\u041e\u043d\u0438 \u0438\u0441\u043f\u043e\u043b\u044c\u0437\u0443\u044e\u0438\u0441\u044f \u0432 \u043d\u0430\u043f\u0438\u0441\u0430\u043d\u0438\u0438 \u0441\u0446\u0435\u043d\u0430\u0440\u0438\u0435\u0432...\" -> \"\u041e\u043d\u0438 \u0438\u0441\u043f\u043e\u043b\u044c\u0437\u0443\u044e\u0442\u0441\u044f \u0432 \u043d\u0430\u043f\u0438\u0441\u0430\u043d\u0438\u0438...\").\n# Word re-seated its hidden \"_GoBack\" (last-edit-position) bookmark from its old spot\n# (an otherwise-empty paragraph further down) to the exact point of this edit, right\n# after the inserted letter. We reproduce both effects below.\n\n$d = $word.ActiveDocument\n\n# 1) Drop the old \"_GoBack\" bookmark wherever it currently lives (Word only ever\n#    keeps one), so it doesn't linger at its old location once we re-add it below.\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks(\"_GoBack\").Delete()\n}\n\n# 2) Locate the misspelled word.\n$findRange = $d.Content\n$findRange.Find.ClearFormatting()\n$findRange.Find.Execute(\"\u0438\u0441\u043f\u043e\u043b\u044c\u0437\u0443\u044e\u0438\u0441\u044f\") | Out-Null\n$typoStart = $findRange.Start\n\n# 3) The fix only ever touches a single character: \"\u0438\u0441\u043f\u043e\u043b\u044c\u0437\u0443\u044e\" + \"\u0438\u0441\u044f\" -> \"\u0438\u0441\u043f\u043e\u043b\u044c\u0437\u0443\u044e\" + \"\u0442\u0441\u044f\",\n#    i.e. the \"\u0438\" right after \"\u0438\u0441\u043f\u043e\u043b\u044c\u0437\u0443\u044e\" (character index 9 of the 12-character typo) becomes \"\u0442\".\n$badChar = $d.Range($typoStart + 9, $typoStart + 10)\n$badChar.Text = \"\u0442\"\n\n# Nudge the formatting so the host keeps this \"\u0442\" as its own run instead of silently\n# re-merging it with its neighbor - matching the author's document, where it stayed a\n# distinct run next to the new bookmark.\n$origBold = $badChar.Font.Bold\n$badChar.Font.Bold = $true\n$badChar.Font.Bold = $origBold\n\n# 4) Re-insert \"_GoBack\" as a zero-length bookmark immediately after the new \"\u0442\",\n#    i.e. right before \"\u0441\u044f ...\" - exactly where Word left the cursor after this edit.\n$newGoBackPoint = $d.Range($typoStart + 10, $typoStart + 10)\n$d.Bookmarks.Add(\"_GoBack\", $newGoBackPoint) | Out-Null\n"}
